$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamps = @(
    45995.01041666666,
    45995.02083333334,
    45995.03125,
    45995.04166666666,
    45995.05208333334,
    45995.0625,
    45995.07291666666,
    45995.08333333334,
    45995.09375,
    45995.10416666666,
    45995.11458333334,
    45995.125,
    45995.13541666666,
    45995.14583333334,
    45995.15625,
    45995.16666666666,
    45995.17708333334,
    45995.1875,
    45995.19791666666,
    45995.20833333334,
    45995.21875,
    45995.22916666666,
    45995.23958333334,
    45995.25,
    45995.26041666666,
    45995.27083333334,
    45995.28125,
    45995.29166666666,
    45995.30208333334,
    45995.3125,
    45995.32291666666,
    45995.33333333334,
    45995.34375,
    45995.35416666666,
    45995.36458333334,
    45995.375,
    45995.38541666666,
    45995.39583333334,
    45995.40625,
    45995.41666666666,
    45995.42708333334,
    45995.4375,
    45995.44791666666,
    45995.45833333334,
    45995.46875,
    45995.47916666666,
    45995.48958333334,
    45995.5,
    45995.51041666666,
    45995.52083333334,
    45995.53125,
    45995.54166666666,
    45995.55208333334,
    45995.5625,
    45995.57291666666,
    45995.58333333334,
    45995.59375,
    45995.60416666666,
    45995.61458333334,
    45995.625,
    45995.63541666666,
    45995.64583333334,
    45995.65625,
    45995.66666666666,
    45995.67708333334,
    45995.6875,
    45995.69791666666,
    45995.70833333334,
    45995.71875,
    45995.72916666666,
    45995.73958333334,
    45995.75,
    45995.76041666666,
    45995.77083333334,
    45995.78125,
    45995.79166666666,
    45995.80208333334,
    45995.8125,
    45995.82291666666,
    45995.83333333334,
    45995.84375,
    45995.85416666666,
    45995.86458333334,
    45995.875,
    45995.88541666666,
    45995.89583333334,
    45995.90625,
    45995.91666666666,
    45995.92708333334,
    45995.9375,
    45995.94791666666,
    45995.95833333334,
    45995.96875,
    45995.97916666666,
    45995.98958333334,
    45996
)

$values = @(
    0.442,
    0.434,
    0,
    0.426,
    0.518,
    0.51,
    0.506,
    0,
    0.51,
    0.494,
    0.506,
    0.514,
    0.598,
    0.594,
    0,
    0,
    5.502,
    0,
    5.518,
    5.51,
    9.454000000000001,
    9.33,
    9.406000000000001,
    9.686,
    13.561,
    14.141,
    16.303,
    23.459,
    62.994,
    84,
    114.256,
    150.072,
    295.309,
    339.84,
    387.601,
    433.256,
    569.638,
    612.753,
    656.225,
    696.862,
    778.03,
    804.7380000000001,
    829.624,
    845.447,
    860.853,
    858.6950000000001,
    863.399,
    854.7190000000001,
    824.623,
    799.261,
    765.288,
    725.525,
    606.96,
    556.3049999999999,
    495.643,
    445.952,
    282.983,
    232.576,
    181.091,
    140.724,
    51.331,
    29.849,
    21.251,
    16,
    9.538,
    9.513999999999999,
    7.013,
    7.1,
    1.222,
    1.214,
    1.218,
    0.694,
    0.678,
    0,
    0.6820000000000001,
    0.658,
    0.738,
    0.746,
    0.754,
    0.738,
    0.75,
    0,
    0,
    0.734,
    0.73,
    0.734,
    0.75,
    0.698,
    0.59,
    0.61,
    0.622,
    0.598,
    0,
    0,
    0,
    0
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $timestamps[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}